$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 289
$ws1.Range("F5").Value = 157
$ws1.Range("F6").Value = 103
$ws1.Range("F7").Value = 285
$ws1.Range("F9").Value = 2027
$ws1.Range("F11").Value = 4815

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G3").Value = 135

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G5").Value = 135
$ws4.Range("F6").Value = 289
$ws4.Range("F7").Value = 157
$ws4.Range("F8").Value = 103
$ws4.Range("F9").Value = 285
$ws4.Range("F13").Value = 2027
$ws4.Range("F15").Value = 4815
